# Update the "想去人数" (number of people interested) figures that were
# refreshed when gh-pages output was regenerated at commit 456a3b4.
#
# Affected sheets: 展览 and 全部类型 (演出 and 本地生活 are unchanged).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F4").Value = 3684
    $ws.Range("F8").Value = 7

    if ($name -eq "展览") {
        $ws.Range("F11").Value = 78
        $ws.Range("F14").Value = 2091
        $ws.Range("F15").Value = 151
    } elseif ($name -eq "全部类型") {
        $ws.Range("F12").Value = 78
        $ws.Range("F17").Value = 2091
        $ws.Range("F18").Value = 151
    }
}
